$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 288.4
$ws.Range("I55").Value = 324.25
$ws.Range("J55").Value = 145
$ws.Range("K55").Value = 324.25
$ws.Range("L55").Value = 145
$ws.Range("M55").Value = -110.25
$ws.Range("N55").Value = -573

$ws.Range("H58").Value = 8531.277
$ws.Range("I58").Value = 357
$ws.Range("J58").Value = 24879.834
$ws.Range("K58").Value = 1071
$ws.Range("L58").Value = 74639.50199999999
$ws.Range("M58").Value = -921
$ws.Range("N58").Value = -74939.50199999999

$ws.Range("H70").Value = 1141.4584
$ws.Range("I70").Value = 1018.82355
$ws.Range("J70").Value = 1439.2858
$ws.Range("K70").Value = 3056.47065
$ws.Range("L70").Value = 4317.857400000001
$ws.Range("M70").Value = -2786.47065
$ws.Range("N70").Value = -4857.857400000001

$ws.Range("H73").Value = 1141.4584
$ws.Range("I73").Value = 1018.82355
$ws.Range("J73").Value = 1439.2858
$ws.Range("K73").Value = 3056.47065
$ws.Range("L73").Value = 4317.857400000001
$ws.Range("M73").Value = -2120.47065
$ws.Range("N73").Value = -6189.857400000001

$ws.Range("H103").Value = 508.66666
$ws.Range("J103").Value = 527
$ws.Range("L103").Value = 1581
$ws.Range("N103").Value = -2753

$ws.Range("H131").Value = 2600.375
$ws.Range("I131").Value = 2600.375
$ws.Range("K131").Value = 7801.125
$ws.Range("M131").Value = -2761.125

$ws.Range("H138").Value = 83336960

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 40034
$ws.Range("I37").Value = 40034
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 40034
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -39761

$ws.Range("H132").Value = 11117.979
$ws.Range("I132").Value = 6958.4243
$ws.Range("J132").Value = 21676.846
$ws.Range("K132").Value = 20875.2729
$ws.Range("L132").Value = 65030.538
$ws.Range("M132").Value = -18345.2729
$ws.Range("N132").Value = -70090.538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1284.5962
$ws.Range("I94").Value = 920.5769
$ws.Range("J94").Value = 1648.6154
$ws.Range("K94").Value = 920.5769
$ws.Range("L94").Value = 1648.6154
$ws.Range("M94").Value = -469.5769
$ws.Range("N94").Value = -2550.6154

$ws.Range("H105").Value = 3199.0454
$ws.Range("I105").Value = 3201.25
$ws.Range("K105").Value = 3201.25
$ws.Range("M105").Value = -1454.25

$ws.Range("H107").Value = 677.5833
$ws.Range("I107").Value = 610
$ws.Range("K107").Value = 610
$ws.Range("M107").Value = 1310

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1646.0454
$ws.Range("I31").Value = 1422.3572
$ws.Range("J31").Value = 2037.5
$ws.Range("K31").Value = 1422.3572
$ws.Range("L31").Value = 2037.5
$ws.Range("M31").Value = -1127.3572
$ws.Range("N31").Value = -2627.5

$ws.Range("H34").Value = 1646.0454
$ws.Range("I34").Value = 1422.3572
$ws.Range("J34").Value = 2037.5
$ws.Range("K34").Value = 1422.3572
$ws.Range("L34").Value = 2037.5
$ws.Range("M34").Value = -1220.3572
$ws.Range("N34").Value = -2441.5

$ws.Range("H39").Value = 6475
$ws.Range("I39").Value = 6475
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 6475
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -6084

$ws.Range("H49").Value = 6475
$ws.Range("I49").Value = 6475
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 6475
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -6293

$ws.Range("H122").Value = 2746.45
$ws.Range("I122").Value = 2296.6667
$ws.Range("J122").Value = 3114.4546
$ws.Range("K122").Value = 6890.000100000001
$ws.Range("L122").Value = 9343.363799999999
$ws.Range("M122").Value = -4440.000100000001
$ws.Range("N122").Value = -14243.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 71428740
$ws.Range("I4").Value = 76923256
$ws.Range("K4").Value = 230769768
$ws.Range("M4").Value = -230769656

$ws.Range("H14").Value = 589.1818
$ws.Range("I14").Value = 589.1818
$ws.Range("K14").Value = 1767.5454
$ws.Range("M14").Value = -1594.5454

$ws.Range("H37").Value = 112142.14
$ws.Range("J37").Value = 112142.14
$ws.Range("L37").Value = 336426.42
$ws.Range("N37").Value = -336650.42

$ws.Range("H62").Value = 9342.857
$ws.Range("I62").Value = 6466.6665
$ws.Range("K62").Value = 19399.9995
$ws.Range("M62").Value = -18713.9995

$ws.Range("H65").Value = 9342.857
$ws.Range("I65").Value = 6466.6665
$ws.Range("K65").Value = 58199.9985
$ws.Range("M65").Value = -54767.9985

$ws.Range("H92").Value = 208.81818
$ws.Range("J92").Value = 124.28571
$ws.Range("L92").Value = 372.85713
$ws.Range("N92").Value = -2868.85713

$ws.Range("H121").Value = 4235.3335
$ws.Range("J121").Value = 4524.909
$ws.Range("L121").Value = 13574.727
$ws.Range("N121").Value = -16194.727

$ws.Range("H129").Value = 4466.2
$ws.Range("J129").Value = 2048
$ws.Range("L129").Value = 6144
$ws.Range("N129").Value = -16144

$ws.Range("H131").Value = 1960.6459
$ws.Range("I131").Value = 1602.6666
$ws.Range("J131").Value = 2011.7858
$ws.Range("K131").Value = 4807.9998
$ws.Range("L131").Value = 6035.357400000001
$ws.Range("M131").Value = 232.0002000000004
$ws.Range("N131").Value = -16115.3574

$ws.Range("H132").Value = 1505.8
$ws.Range("I132").Value = 1429.8462
$ws.Range("J132").Value = 1999.5
$ws.Range("K132").Value = 12868.6158
$ws.Range("L132").Value = 17995.5
$ws.Range("M132").Value = -10338.6158
$ws.Range("N132").Value = -23055.5

$ws.Range("H139").Value = 63718.4
$ws.Range("I139").Value = 79148.5
$ws.Range("J139").Value = 1998
$ws.Range("K139").Value = 237445.5
$ws.Range("L139").Value = 5994
$ws.Range("M139").Value = -232305.5
$ws.Range("N139").Value = -16274

$ws.Range("H140").Value = 1248.7391
$ws.Range("I140").Value = 891.4761999999999
$ws.Range("K140").Value = 2674.4286
$ws.Range("M140").Value = 2505.5714

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 2477.5
$ws.Range("I10").Value = 2455
$ws.Range("J10").Value = 2500
$ws.Range("K10").Value = 2455
$ws.Range("L10").Value = 2500
$ws.Range("M10").Value = -2286
$ws.Range("N10").Value = -2838

$ws.Range("H102").Value = 2431.652
$ws.Range("I102").Value = 1753.3529
$ws.Range("K102").Value = 1753.3529
$ws.Range("M102").Value = -131.3529000000001

$ws.Range("H107").Value = 1510.8
$ws.Range("I107").Value = 2632.1667
$ws.Range("K107").Value = 2632.1667
$ws.Range("M107").Value = -712.1667000000002

$ws.Range("H132").Value = 3912.75
$ws.Range("I132").Value = 3503.6667
$ws.Range("J132").Value = 5140
$ws.Range("K132").Value = 10511.0001
$ws.Range("L132").Value = 15420
$ws.Range("M132").Value = -7981.000100000001
$ws.Range("N132").Value = -20480

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6299.1333
$ws.Range("J7").Value = 9063.125
$ws.Range("L7").Value = 9063.125
$ws.Range("N7").Value = -9287.125

$ws.Range("H126").Value = 6299.1333
$ws.Range("J126").Value = 9063.125
$ws.Range("L126").Value = 27189.375
$ws.Range("N126").Value = -32129.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 10421261
$ws.Range("I122").Value = 15629442
$ws.Range("J122").Value = 4900.875
$ws.Range("K122").Value = 46888326
$ws.Range("L122").Value = 14702.625
$ws.Range("M122").Value = -46885876
$ws.Range("N122").Value = -19602.625

$ws.Range("H126").Value = 15153064
$ws.Range("I126").Value = 18520140
$ws.Range("J126").Value = 1225.25
$ws.Range("K126").Value = 55560420
$ws.Range("L126").Value = 3675.75
$ws.Range("M126").Value = -55557950
$ws.Range("N126").Value = -8615.75
